$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates; force text number-format while writing so Excel's
# COM auto-type-detection doesn't coerce numeric-looking strings (e.g.
# "266.90", "0.0858") into floating point numbers and lose formatting,
# then restore the Normal style so no stray number format lingers.
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "51.511.34"
Set-TextValue "E2" "  -0.50%  "
Set-TextValue "D3" "3.102.98"
Set-TextValue "E3" "  +2.20%  "
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "D5" "386.05"
Set-TextValue "E5" "  +1.22%  "
Set-TextValue "D6" "104.09"
Set-TextValue "E6" "  +0.55%  "
Set-TextValue "E7" "  -1.30%  "
Set-TextValue "E8" "  +0.01%  "
Set-TextValue "D9" "0.588"
Set-TextValue "E9" "  -1.64%  "
Set-TextValue "D10" "37.25"
Set-TextValue "E10" "  +0.09%  "
Set-TextValue "E11" "  +0.06%  "
Set-TextValue "D12" "0.0858"
Set-TextValue "E12" "  -0.45%  "
Set-TextValue "D13" "3.594.06"
Set-TextValue "E13" "  +1.95%  "
Set-TextValue "D14" "18.59"
Set-TextValue "E14" "  -0.23%  "
Set-TextValue "E15" "  +1.01%  "
Set-TextValue "D16" "3.094.16"
Set-TextValue "E16" "  +1.20%  "
Set-TextValue "E17" "  +1.83%  "
Set-TextValue "D18" "10.95"
Set-TextValue "E18" "  +3.79%  "
Set-TextValue "D19" "51.585.38"
Set-TextValue "E19" "  -0.30%  "
Set-TextValue "E20" "  +7.16%  "
Set-TextValue "D21" "12.55"
Set-TextValue "E21" "  -0.24%  "
Set-TextValue "D22" "0.0₃0966"
Set-TextValue "E22" "  +0.00%  "
Set-TextValue "E23" "  +0.02%  "
Set-TextValue "D24" "266.90"
Set-TextValue "E24" "  -0.90%  "
Set-TextValue "D25" "3.18"
Set-TextValue "E25" "  +0.20%  "
Set-TextValue "D26" "8.12"
Set-TextValue "E26" "  -0.81%  "
Set-TextValue "D27" "27.43"
Set-TextValue "E27" "  +4.03%  "
Set-TextValue "D28" "7.23"
Set-TextValue "E28" "  -4.29%  "
Set-TextValue "E29" "  +0.03%  "
Set-TextValue "E30" "  -3.95%  "
Set-TextValue "D31" "0.107"
Set-TextValue "E31" "  -2.18%  "
Set-TextValue "E32" "  +0.94%  "
Set-TextValue "B33" "InjectiveProtocol"
Set-TextValue "C33" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D33" "35.70"
Set-TextValue "E33" "  +3.96%  "
Set-TextValue "B34" "VeChain"
Set-TextValue "C34" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D34" "0.0478"
Set-TextValue "E34" "  +5.50%  "
Set-TextValue "E35" "  +0.86%  "
Set-TextValue "D36" "50.08"
Set-TextValue "E36" "  -0.80%  "
Set-TextValue "E37" "  -0.09%  "
Set-TextValue "D38" "3.39"
Set-TextValue "E38" "  +0.96%  "
Set-TextValue "E39" "  +0.86%  "
Set-TextValue "E40" "  -0.15%  "
Set-TextValue "D41" "128.99"
Set-TextValue "E41" "  +1.29%  "
Set-TextValue "D42" "16.67"
Set-TextValue "E42" "  -2.89%  "
Set-TextValue "E43" "  -0.44%  "
Set-TextValue "E44" "  -2.88%  "
Set-TextValue "D45" "3.78"
Set-TextValue "E45" "  +0.40%  "
Set-TextValue "D46" "22.26"
Set-TextValue "E46" "  +1.75%  "
Set-TextValue "E47" "  +4.80%  "
Set-TextValue "E48" "  -2.57%  "
Set-TextValue "D49" "2.077.62"
Set-TextValue "E49" "  +1.88%  "
Set-TextValue "D50" "0.0336"
Set-TextValue "E50" "  +4.84%  "
